$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and name it "importFrom-ape"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "importFrom-ape"

# List the functions imported from the `ape` package
$values = @("read.tree", "write.tree", "root", "unroot", "is.binary.tree", "multi2di")
for ($i = 0; $i -lt $values.Count; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Make the new sheet the active sheet / tab, matching the selection in the source file
$newSheet.Activate()
$newSheet.Range("A6").Select()
